# ---------------------------------------------------------------------------
# "19. Regras de negócio.docx" — correção de todos os artefatos
#
# 1. Center the title paragraph.
# 2. RN0001-RN0006 paragraphs: merge the "RNxxxx" + ": " runs into a single
#    bold "RNxxxx: " run, and replace the body text (re-shuffled / fixed
#    business-rule sentences).
# 3. RN0007 / RN0008 paragraphs: keep the bold "RNxxxx" label run as-is and
#    just replace the body run's text (which still carries the leading
#    ": ").
# 4. Remove the RN0009 and RN0010 paragraphs entirely.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Center the title paragraph -----------------------------------------
$d.Paragraphs.Item(1).Alignment = 1   # wdAlignParagraphCenter

# --- 2. RN0001-RN0006: merge label+colon run, replace body -----------------
# Each entry: paragraph index -> new full paragraph text ("RNxxxx: body").
# The bold portion is always the "RNxxxx: " prefix (8 characters).

function Set-RuleParagraph($paraIndex, $fullText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $pEndNoMark = $p.Range.End - 1

    $textRng = $d.Range($pStart, $pEndNoMark)
    $textRng.Text = $fullText

    $labelLen = 8   # length of "RNxxxx: "
    $bodyStart = $pStart + $labelLen
    $bodyEnd = $pStart + $fullText.Length
    if ($bodyEnd -gt $bodyStart) {
        $bodyRng = $d.Range($bodyStart, $bodyEnd)
        $bodyRng.Font.Bold = 0
    }
}

Set-RuleParagraph 3 'RN0001: Somente serão aceitos pedidos ao qual o cliente informa as características do produto desejado.'
Set-RuleParagraph 4 'RN0002: Clientes que realizarem o pedido no qual o valor total seja acima de R$100,00 terão frete grátis.'
Set-RuleParagraph 5 'RN0003: O pedido será processado somente quando o cliente informar todos os dados obrigatórios.'
Set-RuleParagraph 6 'RN0004: O produto poderá ser devolvido se atender às características de garantia e/ou produto com defeito.'
Set-RuleParagraph 7 'RN0005: O produto pode ser cancelado antes da entrega, atendendo a política de cancelamento da loja.'
Set-RuleParagraph 8 'RN0006: Para a confirmação dos dados do cliente, será feito uma autenticação por e-mail e/ou celular.'

# --- 3. RN0007 / RN0008: keep label run, replace body run only -------------

function Set-RuleBodyOnly($paraIndex, $newBodyText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $pEndNoMark = $p.Range.End - 1

    $labelLen = 6   # length of "RNxxxx" (no colon, kept in the body run)
    $bodyRng = $d.Range($pStart + $labelLen, $pEndNoMark)
    $bodyRng.Text = $newBodyText
}

Set-RuleBodyOnly 9 ': Para gerar a nota fiscal da compra, é necessário a confirmação do pagamento.'
Set-RuleBodyOnly 10 ': Os valores dos produtos apenas serão alterados pelo proprietário.'

# --- 4. Remove the RN0009 and RN0010 paragraphs -----------------------------
# Delete from the highest index down so indices of earlier paragraphs stay
# valid.
$d.Paragraphs.Item(12).Range.Delete()
$d.Paragraphs.Item(11).Range.Delete()
